$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 283, pushing the current
# rows 283:296 down to 285:298 (matches the dimension growing from
# A1:R296 to A1:R298).
$ws.Rows.Item(283).Insert()
$ws.Rows.Item(283).Insert()

# New row 283
$ws.Cells.Item(283, 1).Value = 10
$ws.Cells.Item(283, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(283, 3).Value = "La Araucanía"
$ws.Cells.Item(283, 4).Value = 44585
$ws.Cells.Item(283, 5).Value = 9
$ws.Cells.Item(283, 6).Value = 100112037
$ws.Cells.Item(283, 7).Value = "Cebollín"
$ws.Cells.Item(283, 8).Value = "Sin especificar"
$ws.Cells.Item(283, 9).Value = "Primera"
$ws.Cells.Item(283, 10).Value = 55
$ws.Cells.Item(283, 11).Value = 7000
$ws.Cells.Item(283, 12).Value = 7000
$ws.Cells.Item(283, 13).Value = 7000
$ws.Cells.Item(283, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(283, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(283, 16).Value = 583
$ws.Cells.Item(283, 17).Value = 12
$ws.Cells.Item(283, 18).Value = "Hortaliza"

# New row 284
$ws.Cells.Item(284, 1).Value = 10
$ws.Cells.Item(284, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(284, 3).Value = "La Araucanía"
$ws.Cells.Item(284, 4).Value = 44585
$ws.Cells.Item(284, 5).Value = 9
$ws.Cells.Item(284, 6).Value = 100112037
$ws.Cells.Item(284, 7).Value = "Cebollín"
$ws.Cells.Item(284, 8).Value = "Sin especificar"
$ws.Cells.Item(284, 9).Value = "Primera"
$ws.Cells.Item(284, 10).Value = 110
$ws.Cells.Item(284, 11).Value = 5000
$ws.Cells.Item(284, 12).Value = 5000
$ws.Cells.Item(284, 13).Value = 5000
$ws.Cells.Item(284, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(284, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(284, 16).Value = 417
$ws.Cells.Item(284, 17).Value = 12
$ws.Cells.Item(284, 18).Value = "Hortaliza"
